$d = $word.ActiveDocument
$t = $d.Tables(1)

# Rows 1-12: summary statistics column, update single-value cells in place.
$t.Rows(1).Cells(1).Range.Text = "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"
$t.Rows(4).Cells(1).Range.Text = "269"
$t.Rows(5).Cells(1).Range.Text = "0.00003"
$t.Rows(6).Cells(1).Range.Text = "0.00282"
$t.Rows(7).Cells(1).Range.Text = "0.00019"
$t.Rows(8).Cells(1).Range.Text = "0.00011"
$t.Rows(9).Cells(1).Range.Text = "0.00027"
$t.Rows(10).Cells(1).Range.Text = "0.00033"
$t.Rows(11).Cells(1).Range.Text = "0.00043"
$t.Rows(12).Cells(1).Range.Text = "0.06081"

# Rows 44-46: collapse the tab-separated per-iteration detail rows down to
# their single aggregate value (matching rows 1-3 above).
$t.Rows(44).Cells(1).Range.Text = "100"
$t.Rows(45).Cells(1).Range.Text = "0.06"
$t.Rows(46).Cells(1).Range.Text = "2432"
